$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 909.0625
$ws.Range("J17").Value = 909.0625
$ws.Range("L17").Value = 2727.1875
$ws.Range("N17").Value = -3063.1875
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1000
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -1350
$ws.Range("H132").Value = 733440.06
$ws.Range("I132").Value = 2092.7917
$ws.Range("J132").Value = 2581054.2
$ws.Range("K132").Value = 6278.375100000001
$ws.Range("L132").Value = 7743162.600000001
$ws.Range("M132").Value = -3748.375100000001
$ws.Range("N132").Value = -7748222.600000001
$ws.Range("H138").Value = 2224288.8
$ws.Range("I138").Value = 1308.0613
$ws.Range("J138").Value = 6413752.5
$ws.Range("K138").Value = 3924.1839
$ws.Range("L138").Value = 19241257.5
$ws.Range("M138").Value = 1215.8161
$ws.Range("N138").Value = -19251537.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 40382.055
$ws.Range("I132").Value = 29269.97
$ws.Range("J132").Value = 61988.89
$ws.Range("K132").Value = 87809.91
$ws.Range("L132").Value = 185966.67
$ws.Range("M132").Value = -85279.91
$ws.Range("N132").Value = -191026.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H28").Value = 21500
$ws.Range("J28").Value = 21500
$ws.Range("L28").Value = 21500
$ws.Range("N28").Value = -22088
$ws.Range("H134").Value = 1856.8937
$ws.Range("I134").Value = 985.96875
$ws.Range("J134").Value = 3714.8667
$ws.Range("K134").Value = 2957.90625
$ws.Range("L134").Value = 11144.6001
$ws.Range("M134").Value = -422.90625
$ws.Range("N134").Value = -16214.6001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 14929.603
$ws.Range("I132").Value = 1046.125
$ws.Range("J132").Value = 60663.41
$ws.Range("K132").Value = 3138.375
$ws.Range("L132").Value = 181990.23
$ws.Range("M132").Value = -608.375
$ws.Range("N132").Value = -187050.23
$ws.Range("H134").Value = 15929.575
$ws.Range("I134").Value = 1115.9623
$ws.Range("J134").Value = 55185.65
$ws.Range("K134").Value = 3347.8869
$ws.Range("L134").Value = 165556.95
$ws.Range("M134").Value = -812.8868999999995
$ws.Range("N134").Value = -170626.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 237
$ws.Range("I18").Value = 121.25
$ws.Range("J18").Value = 700
$ws.Range("K18").Value = 363.75
$ws.Range("L18").Value = 2100
$ws.Range("M18").Value = -194.75
$ws.Range("N18").Value = -2438
$ws.Range("H109").Value = 3041.8
$ws.Range("I109").Value = 709
$ws.Range("J109").Value = 3625
$ws.Range("K109").Value = 2127
$ws.Range("L109").Value = 10875
$ws.Range("M109").Value = -1087
$ws.Range("N109").Value = -12955
$ws.Range("H111").Value = 2563.5
$ws.Range("I111").Value = 418
$ws.Range("J111").Value = 9000
$ws.Range("K111").Value = 1254
$ws.Range("L111").Value = 27000
$ws.Range("M111").Value = 1813
$ws.Range("N111").Value = -33134
$ws.Range("H112").Value = 17547414
$ws.Range("I112").Value = 2305.4
$ws.Range("J112").Value = 23813524
$ws.Range("K112").Value = 6916.200000000001
$ws.Range("L112").Value = 71440572
$ws.Range("M112").Value = -5808.200000000001
$ws.Range("N112").Value = -71442788
$ws.Range("H113").Value = 442.125
$ws.Range("I113").Value = 330.23077
$ws.Range("K113").Value = 990.69231
$ws.Range("M113").Value = 1179.30769
$ws.Range("H115").Value = 2380.0476
$ws.Range("J115").Value = 2349.05
$ws.Range("L115").Value = 7047.150000000001
$ws.Range("N115").Value = -9397.150000000001
$ws.Range("H116").Value = 112604.25
$ws.Range("I116").Value = 222696.33
$ws.Range("J116").Value = 75906.89
$ws.Range("K116").Value = 668088.99
$ws.Range("L116").Value = 227720.67
$ws.Range("M116").Value = -664646.99
$ws.Range("N116").Value = -234604.67
$ws.Range("H122").Value = 690.44446
$ws.Range("I122").Value = 288.33334
$ws.Range("J122").Value = 977.6667
$ws.Range("K122").Value = 2595.00006
$ws.Range("L122").Value = 8799.0003
$ws.Range("M122").Value = -145.0000600000003
$ws.Range("N122").Value = -13699.0003
$ws.Range("H127").Value = 1883.25
$ws.Range("J127").Value = 1883.25
$ws.Range("L127").Value = 5649.75
$ws.Range("N127").Value = -15569.75
$ws.Range("H130").Value = 2924
$ws.Range("J130").Value = 3148.889
$ws.Range("L130").Value = 9446.667000000001
$ws.Range("N130").Value = -19486.667
$ws.Range("H131").Value = 1241.5
$ws.Range("J131").Value = 1461.0416
$ws.Range("L131").Value = 4383.1248
$ws.Range("N131").Value = -14463.1248

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 15119
$ws.Range("J136").Value = 15119
$ws.Range("L136").Value = 45357
$ws.Range("N136").Value = -50457

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 719.91895
$ws.Range("I22").Value = 470.83334
$ws.Range("J22").Value = 839.48
$ws.Range("K22").Value = 470.83334
$ws.Range("L22").Value = 839.48
$ws.Range("M22").Value = -175.83334
$ws.Range("N22").Value = -1429.48
$ws.Range("H27").Value = 719.91895
$ws.Range("I27").Value = 470.83334
$ws.Range("J27").Value = 839.48
$ws.Range("K27").Value = 470.83334
$ws.Range("L27").Value = 839.48
$ws.Range("M27").Value = -363.83334
$ws.Range("N27").Value = -1053.48
$ws.Range("H132").Value = 21967.49
$ws.Range("I132").Value = 1097.303
$ws.Range("J132").Value = 65012.25
$ws.Range("K132").Value = 3291.909000000001
$ws.Range("L132").Value = 195036.75
$ws.Range("M132").Value = -761.9090000000006
$ws.Range("N132").Value = -200096.75
$ws.Range("H136").Value = 35478.69
$ws.Range("I136").Value = 22035.852
$ws.Range("J136").Value = 92916.27
$ws.Range("K136").Value = 66107.556
$ws.Range("L136").Value = 278748.81
$ws.Range("M136").Value = -63557.556
$ws.Range("N136").Value = -283848.81

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 27480.28
$ws.Range("I132").Value = 16940.459
$ws.Range("J132").Value = 73403.79
$ws.Range("K132").Value = 50821.37699999999
$ws.Range("L132").Value = 220211.37
$ws.Range("M132").Value = -48291.37699999999
$ws.Range("N132").Value = -225271.37
$ws.Range("H136").Value = 36888.492
$ws.Range("I136").Value = 27855.244
$ws.Range("J136").Value = 53600
$ws.Range("K136").Value = 83565.73199999999
$ws.Range("L136").Value = 160800
$ws.Range("M136").Value = -81015.73199999999
$ws.Range("N136").Value = -165900

Write-Output "done"